# Update cryptos list figures (price + 1h volume change) and fix the
# row ordering for a few coins (Toncoin/Dogecoin, Mantle/Filecoin,
# EnergySwap/dogwifhat/Cosmos) to match the latest scrape.
#
# Numeric-looking price strings (e.g. "1.00", "597.64") are written with a
# leading apostrophe to force Excel to keep them as text (matching the
# original inlineStr cells) instead of auto-converting them to numbers;
# the style is then reset to "Normal" so the quote-prefix formatting isn't
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.401.56'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '3.500.87'
$ws.Range("E3").Value = '  -2.49%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = "'597.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = "'141.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.85%  '
$ws.Range("D7").Value = '3.498.85'
$ws.Range("E7").Value = '  -2.53%  '
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = "'0.512"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.34%  '
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = "'7.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.80%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = "'0.130"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.85%  '
$ws.Range("D12").Value = "'0.401"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.84%  '
$ws.Range("D13").Value = '4.111.34'
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").Value = "'0.0000193"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.83%  '
$ws.Range("D15").Value = "'28.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.01%  '
$ws.Range("D16").Value = '3.509.39'
$ws.Range("E16").Value = '  -2.12%  '
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").Value = '65.513.81'
$ws.Range("E18").Value = '  -1.85%  '
$ws.Range("D19").Value = "'10.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.47%  '
$ws.Range("D20").Value = "'6.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.61%  '
$ws.Range("D21").Value = "'14.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.83%  '
$ws.Range("D22").Value = "'415.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.77%  '
$ws.Range("D23").Value = "'0.590"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.77%  '
$ws.Range("D24").Value = "'76.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.87%  '
$ws.Range("D25").Value = '3.655.26'
$ws.Range("E25").Value = '  -2.16%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = "'0.0000112"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.73%  '
$ws.Range("D28").Value = "'2.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.60%  '
$ws.Range("D29").Value = "'7.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.10%  '
$ws.Range("D30").Value = "'8.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.99%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '3.522.68'
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("D33").Value = "'0.152"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("D34").Value = "'24.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.88%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -8.33%  '
$ws.Range("D37").Value = "'7.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.40%  '
$ws.Range("D38").Value = "'174.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'5.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.79%  '
$ws.Range("D40").Value = "'1.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.05%  '
$ws.Range("D41").Value = "'0.0804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.38%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = "'0.850"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.29%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = "'4.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.68%  '
$ws.Range("D44").Value = "'45.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.25%  '
$ws.Range("D45").Value = "'1.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.94%  '
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'23.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.14%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = "'2.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.65%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = "'6.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("D50").Value = "'1.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.41%  '
$ws.Range("D51").Value = "'0.895"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.99%  '
